$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new "SKU Code" column in front of the sales figures -
#    this shifts the existing "T-x month ..." columns one to the
#    right and leaves a blank column B to fill with the id values.
#    Give the new header the same look as the other headers.
# ------------------------------------------------------------------
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("C1").Copy($ws.Range("B1"))
$ws.Range("B1").Value = "SKU Code"
$excel.CutCopyMode = 0

# the SKU values used to live in (what is now) column C - move them
# into the new column B and clean up behind them.
$ws.Range("C2:C3").Copy($ws.Range("B2"))
$ws.Range("C2:C3").ClearContents()
$excel.CutCopyMode = 0

# it's an identifier, not a quantity - drop the thousands-style
# numeric formatting that the sales columns use, keep the same
# look otherwise.
$ws.Range("B2:B3").ClearFormats()
$ws.Range("B2:B3").Font.Name = "Calibri"
$ws.Range("B2:B3").Font.Size = 11
$ws.Range("B2:B3").Font.ColorIndex = 1
$ws.Range("B2:B3").Interior.Pattern = 1
$ws.Range("B2:B3").Interior.Color = 16777215
$ws.Range("B2:B3").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 2. The "T month" figures (now in I:K) are the latest numbers we
#    have, so carry them across all three reporting periods
#    (T-2, T-1, T) for every channel, replacing the stale
#    per-period figures. Stage them in a scratch range first so the
#    fan-out (which overlaps the source columns) doesn't clobber
#    itself.
# ------------------------------------------------------------------
$ws.Range("I2:K3").Copy()
$ws.Range("Z2").PasteSpecial(-4163)

$ws.Range("Z2:AB3").Copy()
$ws.Range("C2").PasteSpecial(-4163)
$ws.Range("Z2:AB3").Copy()
$ws.Range("F2").PasteSpecial(-4163)
$ws.Range("Z2:AB3").Copy()
$ws.Range("I2").PasteSpecial(-4163)

$ws.Range("Z2:AB3").Clear()
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Drop the trailing "date" / "Month" / "year" columns - no
#    longer needed for this report.
# ------------------------------------------------------------------
$ws.Range("L1:N3").EntireColumn.Delete()

# ------------------------------------------------------------------
# 4. Leave the selection where Excel put it after the edits.
# ------------------------------------------------------------------
$null = $ws.Range("K1").Select()
